$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '26.433.18'
$ws.Range("E2").Value = '  +1.03%  '
$ws.Range("D3").Value = '1.673.91'
$ws.Range("E3").Value = '  +1.14%  '
$ws.Range("E4").Value = '  +0.51%  '
$ws.Range("D5").Value = '''221.74'
$ws.Range("E5").Value = '  +1.77%  '
$ws.Range("D6").Value = '''0.5339'
$ws.Range("E7").Value = '  +0.46%  '
$ws.Range("D8").Value = '''0.2667'
$ws.Range("E8").Value = '  +1.54%  '
$ws.Range("D9").Value = '''0.06396'
$ws.Range("E9").Value = '  +1.27%  '
$ws.Range("D10").Value = '''20.94'
$ws.Range("E10").Value = '  +2.70%  '
$ws.Range("D11").Value = '''0.07856'
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("D12").Value = '''4.535'
$ws.Range("E12").Value = '  +0.47%  '
$ws.Range("D13").Value = '1.680.33'
$ws.Range("E13").Value = '  +1.15%  '
$ws.Range("D14").Value = '1.903.69'
$ws.Range("E14").Value = '  +1.09%  '
$ws.Range("D15").Value = '''0.5629'
$ws.Range("E15").Value = '  +2.60%  '
$ws.Range("D16").Value = '0.0₅8212'
$ws.Range("E16").Value = '  +0.78%  '
$ws.Range("D17").Value = '''66.18'
$ws.Range("E17").Value = '  +1.23%  '
$ws.Range("D18").Value = '26.441.80'
$ws.Range("E18").Value = '  +1.16%  '
$ws.Range("D19").Value = '''1.011'
$ws.Range("E19").Value = '  +0.55%  '
$ws.Range("D20").Value = '''4.729'
$ws.Range("E20").Value = '  +2.95%  '
$ws.Range("D21").Value = '''198.28'
$ws.Range("E21").Value = '  +4.00%  '
$ws.Range("E22").Value = '  +2.60%  '
$ws.Range("E23").Value = '  +1.52%  '
$ws.Range("E24").Value = '  +0.44%  '
$ws.Range("D25").Value = '''146.49'
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("D26").Value = '''0.1229'
$ws.Range("E26").Value = '  +0.50%  '
$ws.Range("D27").Value = '''7.256'
$ws.Range("D28").Value = '''16.28'
$ws.Range("E28").Value = '  +2.04%  '
$ws.Range("D29").Value = '''1.504'
$ws.Range("E29").Value = '  +2.18%  '
$ws.Range("D30").Value = '''0.05932'
$ws.Range("E30").Value = '  +3.74%  '
$ws.Range("E31").Value = '  +1.44%  '
$ws.Range("D32").Value = '''3.564'
$ws.Range("E32").Value = '  +0.48%  '
$ws.Range("D33").Value = '''3.324'
$ws.Range("E33").Value = '  +1.73%  '
$ws.Range("D34").Value = '''1.617'
$ws.Range("E34").Value = '  +1.95%  '
$ws.Range("D35").Value = '''0.9708'
$ws.Range("E35").Value = '  +2.45%  '
$ws.Range("D36").Value = '''2.841'
$ws.Range("E36").Value = '  +1.34%  '
$ws.Range("D37").Value = '''2.438'
$ws.Range("E37").Value = '  +0.67%  '
$ws.Range("D38").Value = '''0.5842'
$ws.Range("E38").Value = '  +2.25%  '
$ws.Range("D39").Value = '''0.01616'
$ws.Range("E39").Value = '  +0.65%  '
$ws.Range("D40").Value = '1.080.01'
$ws.Range("E40").Value = '  +4.00%  '
$ws.Range("D41").Value = '''5.940'
$ws.Range("E41").Value = '  +2.58%  '
$ws.Range("D42").Value = '''0.8671'
$ws.Range("E42").Value = '  +2.04%  '
$ws.Range("D44").Value = '''103.15'
$ws.Range("E44").Value = '  -0.74%  '
$ws.Range("D45").Value = '1.812.84'
$ws.Range("E45").Value = '  +0.98%  '
$ws.Range("D46").Value = '''58.74'
$ws.Range("E46").Value = '  +3.65%  '
$ws.Range("D47").Value = '0.0₈105'
$ws.Range("E47").Value = '  +1.10%  '
$ws.Range("D48").Value = '''1.014'
$ws.Range("E48").Value = '  +1.18%  '
$ws.Range("D49").Value = '''0.4413'
$ws.Range("E49").Value = '  +1.33%  '
$ws.Range("D50").Value = '''8.007'
$ws.Range("E50").Value = '  +1.84%  '
$ws.Range("D51").Value = '''0.05159'
$ws.Range("E51").Value = '  +0.13%  '
